$d = $word.ActiveDocument

# --- Locate the anchor paragraph: the "Man" piece's Interaction bullet -------
# ("Interaction: NONE. The man tells the player the number of hits they can
#  take, and the number of points they need to advance. Has a 1/3 chance of
#  appearing somewhere on the board.") The new "Sage" piece section is added
# directly after it, before the "Friendly Pieces" heading.
$anchorIndex = 0
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Has a 1/3 chance of appearing somewhere on the board.*") {
        $anchorIndex = $p.Index
        break
    }
}

if ($anchorIndex -eq 0) {
    Write-Host "ERROR: could not locate the Man interaction paragraph"
} else {
    # Create one empty paragraph right after the anchor; it inherits the
    # anchor's paragraph formatting (ListParagraph style, numId 1, ilvl 1,
    # Times New Roman run fonts) which we then adjust per new line below.
    $anchor = $d.Paragraphs.Item($anchorIndex)
    $anchor.Range.InsertParagraphAfter()

    # The four new bullet lines that make up the "Sage" piece description,
    # and the outline level each belongs to (1 = top-level bullet like
    # "Man:", 2 = sub-bullet like "Symbol:"/"Motion:"/"Interaction:").
    $sageLines = @(
        "Sage:",
        "Symbol: ‘S’",
        "Motion: Randomly chooses a direction and moves one space.",
        "Interaction: NONE. Randomly informs player about the other piece’s interaction results."
    )
    $sageLevels = @(1, 2, 2, 2)

    for ($i = 0; $i -lt $sageLines.Length; $i++) {
        $curIndex = $anchorIndex + 1 + $i
        $curPara = $d.Paragraphs.Item($curIndex)

        $curPara.Range.InsertAfter($sageLines[$i])

        # Make sure the run carries the same Times New Roman font as the
        # rest of the legend (ascii/hAnsi + complex-script/"cs" slot).
        $curPara.Range.Font.Name = "Times New Roman"
        $curPara.Range.Font.NameBi = "Times New Roman"

        # ListLevelNumber is 1-based (1 => w:ilvl 0, 2 => w:ilvl 1).
        $curPara.Range.ListFormat.ListLevelNumber = $sageLevels[$i]

        if ($i -lt ($sageLines.Length - 1)) {
            $curPara.Range.InsertParagraphAfter()
        }
    }
}
